$d = $word.ActiveDocument

# ---- Change 1: "Votre nom" -> "Vincent Bouchard et Simon Robidas" ----
# The target keeps this as two separate runs ("...Simon R" / "obidas") even
# though both share identical formatting, so once the text is in place we
# briefly drop a bookmark at the seam: saving coalesces same-format runs
# unless something (a bookmark, a field, ...) sits between them, and once
# that temporary bookmark is removed again the seam survives.
$rName = $d.Content
$okName = $rName.Find.Execute("Votre nom", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okName) { throw "Could not find 'Votre nom'" }
$rName.Text = "Vincent Bouchard et Simon Robidas"
$splitPos = $rName.Start + 27
$bmNameRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("ZZZNAME", $bmNameRange) | Out-Null
$d.Bookmarks("ZZZNAME").Delete()

# ---- Change 2: restructure the two paragraphs about the pacman death / title screen ----
# 2a. Drop the trailing period from "... la musique arrete."
$r1 = $d.Content
$ok1 = $r1.Find.Execute("ce que la musique arrête.", $true, $false, $false, $false, $false, $true, 1, $false, "ce que la musique arrête", 2)
if (-not $ok1) { throw "Could not find 'ce que la musique arrête.'" }

# 2b. Append the new clause (plus a throwaway marker char so the insertion
#     point never lands exactly on the paragraph mark - doing that confuses
#     the bookmark placement below).
$r2 = $d.Content
$ok2 = $r2.Find.Execute("ce que la musique arrête", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok2) { throw "Could not find 'ce que la musique arrête'" }
$insPos = $r2.End
$insRange = $d.Range($insPos, $insPos)
$insRange.InsertAfter(" le jeu est en temps d’arrêtZ")

# 2c. Figure out the three run-boundary positions in the new text
$rA = $d.Content
$rA.Find.Execute("jusqu’à", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$posA = $rA.End

$rB = $d.Content
$rB.Find.Execute("ce que la musique arrête", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$posB = $rB.End

$rC = $d.Content
$rC.Find.Execute("le jeu est en temps d’arrêt", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$posC = $rC.End

# 2d. Force the run splits with two throwaway bookmarks, and drop the real
#     _GoBack bookmark at the third seam - this also relocates it away from
#     its old spot between "puis il " and "disparait...".
$bmA = $d.Range($posA, $posA)
$d.Bookmarks.Add("ZZZBMA", $bmA) | Out-Null
$bmB = $d.Range($posB, $posB)
$d.Bookmarks.Add("ZZZBMB", $bmB) | Out-Null
$bmC = $d.Range($posC, $posC)
$d.Bookmarks.Add("_GoBack", $bmC) | Out-Null

$d.Bookmarks("ZZZBMA").Delete()
$d.Bookmarks("ZZZBMB").Delete()

# 2e. Remove the throwaway marker char and put the final period in its place,
#     right after the relocated _GoBack bookmark.
$bmGoBack = $d.Bookmarks("_GoBack")
$zPos = $bmGoBack.End
$zRange = $d.Range($zPos, $zPos + 1)
$zRange.Delete()

$bmGoBack2 = $d.Bookmarks("_GoBack")
$finalPos = $bmGoBack2.End
$finalRange = $d.Range($finalPos, $finalPos)
$finalRange.InsertAfter(".")

Write-Host "All done."
